$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted right before the current row 29,
# pushing the existing rows 29-41 down to rows 30-42.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record's data.
$ws.Range("A29").Value2 = 1
$ws.Range("B29").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C29").Value2 = "Arica y Parinacota"
$ws.Range("D29").Value2 = 44875
$ws.Range("E29").Value2 = 15
$ws.Range("F29").Value2 = 100112045
$ws.Range("G29").Value2 = "Zapallo"
$ws.Range("H29").Value2 = "Camote"
$ws.Range("I29").Value2 = "1a nueva(o)"
$ws.Range("J29").Value2 = 1000
$ws.Range("K29").Value2 = 800
$ws.Range("L29").Value2 = 850
$ws.Range("M29").Value2 = 825
$ws.Range("N29").Value2 = "$/kilo (volumen en unidades)"
$ws.Range("O29").Value2 = "Perú"
$ws.Range("P29").Value2 = 825
$ws.Range("Q29").Value2 = 1
$ws.Range("R29").Value2 = "Hortaliza"
